$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per row (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 5

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 10

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5

$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 10

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 5

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 11

$wb.Save()
